$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - new cells (N2:P2) spelling "SLT"
$ws.Range("N2").Value = "S"
$ws.Range("O2").Value = "L"
$ws.Range("P2").Value = "T"

# Row 15 - new content
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "T"
$ws.Range("C15").Value = "X "
$ws.Range("E15").Value = "t"
$ws.Range("F15").Value = "t"
$ws.Range("G15").Value = "t"
$ws.Range("I15").Value = "b"
$ws.Range("J15").Value = "n"
$ws.Range("K15").Value = "m"
$ws.Range("M15").Value = "p"
$ws.Range("N15").Value = "r"
$ws.Range("O15").Value = "d"
$ws.Range("P15").Value = "B"
$ws.Range("Q15").Value = "m"

# Update selection to N3
$ws.Range("N3").Select()
